# Update the "satrap Trimm Set" (row 3) and "Varta Premium Light F20" (row 4)
# entries by swapping their values (the crawler re-ordered the rows between
# runs), and bump the crawl timestamp in column O for every data row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap the contents of row 3 and row 4 (the crawler re-ordered these two
#     products between runs). Use Range.Copy (not .Value round-tripping) so
#     that cell types (e.g. the text-formatted "44.95" price) survive the
#     swap untouched instead of being re-inferred as numbers. A scratch row
#     well below the used range holds row 3 while row 4 moves into place. ---
$ws.Range("A3:N3").Copy($ws.Range("A100:N100"))
$ws.Range("A4:N4").Copy($ws.Range("A3:N3"))
$ws.Range("A100:N100").Copy($ws.Range("A4:N4"))
$ws.Range("A100:N100").Clear()

# --- Update the crawl timestamp for every data row (2 through 37) ---
for ($r = 2; $r -le 37; $r++) {
    $ws.Range("O$r").Value = "2022-09-13 21:01:02"
}
